$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row 16 that mirrors row 15's layout (new Gaussian Quadrature
# scheme "HexGrid-60degTilt5degRes" entry):
# A16 = 14, formatted like A15 (bold/bordered/centered)
# B16 = same label as B15 ("HexGrid-60degTilt5degRes", shared string)
# C16:M16 = 1

$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A16").Value = 14

$ws.Range("B16").Value = $ws.Range("B15").Value2

$ws.Range("C16:M16").Value = 1
